# [PHOENIX-5917] updated fee details in legacy trade license
#
# The "legencyDetails" sheet holds a row of fee-detail column headers
# (dataName, amount1..amount6, legencyTrade). The header labels for the
# extra fee columns were placeholder letters (amountB, amountC, amountD,
# amountE, amountF) and are renamed to the numbered scheme used elsewhere
# (amount2, amount3, amount4, amount5, amount6) to match amount1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("legencyDetails")

$ws.Range("C1").Value = "amount2"
$ws.Range("D1").Value = "amount3"
$ws.Range("E1").Value = "amount4"
$ws.Range("F1").Value = "amount5"
$ws.Range("G1").Value = "amount6"
